# Swap the contents of columns C (codeforiati:group-name) and D
# (codeforiati:group-code) -- including the header row -- so that the
# group-code column now precedes the group-name column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cValue = $cCell.Value2
    $dValue = $dCell.Value2

    $cCell.Value = $dValue
    $dCell.Value = $cValue
}
